$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.353.21"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3
$ws.Range("D3").Value = "1.877.04"
$ws.Range("E3").Value = "  +0.97%  "

# Row 4
$ws.Range("D4").Value = "'0.9998"

# Row 5
$ws.Range("D5").Value = "'0.7115"
$ws.Range("E5").Value = "  -0.27%  "

# Row 6
$ws.Range("D6").Value = "'242.26"
$ws.Range("E6").Value = "  +0.77%  "

# Row 7
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "'0.3112"
$ws.Range("E8").Value = "  +1.25%  "

# Row 9
$ws.Range("D9").Value = "'0.07767"
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("D10").Value = "'25.10"
$ws.Range("E10").Value = "  +0.76%  "

# Row 11
$ws.Range("D11").Value = "'0.08460"
$ws.Range("E11").Value = "  +2.50%  "

# Row 12
$ws.Range("D12").Value = "1.910.71"
$ws.Range("E12").Value = "  +2.50%  "

# Row 13
$ws.Range("D13").Value = "'5.214"
$ws.Range("E13").Value = "  -0.16%  "

# Row 14
$ws.Range("D14").Value = "'0.7117"
$ws.Range("E14").Value = "  -0.30%  "

# Row 15
$ws.Range("D15").Value = "'91.38"
$ws.Range("E15").Value = "  +1.39%  "

# Row 16
$ws.Range("D16").Value = "29.356.24"
$ws.Range("E16").Value = "  +0.58%  "

# Row 17
$ws.Range("D17").Value = "'0.000008259"
$ws.Range("E17").Value = "  +5.82%  "

# Row 18
$ws.Range("D18").Value = "'6.008"
$ws.Range("E18").Value = "  +2.48%  "

# Row 19
$ws.Range("D19").Value = "'242.53"

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'13.24"
$ws.Range("E20").Value = "  +0.70%  "

# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.126.48"

# Row 22
$ws.Range("D22").Value = "'0.9997"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("D23").Value = "'7.853"
$ws.Range("E23").Value = "  -0.86%  "

# Row 24
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("D25").Value = "'0.1607"
$ws.Range("E25").Value = "  +1.74%  "

# Row 26
$ws.Range("D26").Value = "'162.69"
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("D27").Value = "'9.021"
$ws.Range("E27").Value = "  +1.32%  "

# Row 28
$ws.Range("E28").Value = "  +1.24%  "

# Row 29
$ws.Range("D29").Value = "'1.514"
$ws.Range("E29").Value = "  +1.40%  "

# Row 30
$ws.Range("D30").Value = "'4.405"
$ws.Range("E30").Value = "  +0.51%  "

# Row 31
$ws.Range("D31").Value = "'4.334"
$ws.Range("E31").Value = "  +4.91%  "

# Row 32
$ws.Range("D32").Value = "'1.276"
$ws.Range("E32").Value = "  -3.04%  "

# Row 33
$ws.Range("D33").Value = "'0.05257"
$ws.Range("E33").Value = "  +1.34%  "

# Row 34
$ws.Range("D34").Value = "'1.933"
$ws.Range("E34").Value = "  +1.41%  "

# Row 35
$ws.Range("E35").Value = "  +0.10%  "

# Row 36
$ws.Range("D36").Value = "'0.7419"
$ws.Range("E36").Value = "  +2.03%  "

# Row 37
$ws.Range("D37").Value = "'2.683"
$ws.Range("E37").Value = "  +0.04%  "

# Row 38
$ws.Range("D38").Value = "'0.01869"
$ws.Range("E38").Value = "  +1.29%  "

# Row 39
$ws.Range("D39").Value = "'2.725"
$ws.Range("E39").Value = "  +1.49%  "

# Row 40
$ws.Range("D40").Value = "1.173.40"
$ws.Range("E40").Value = "  +1.59%  "

# Row 41
$ws.Range("D41").Value = "'6.388"
$ws.Range("E41").Value = "  +4.70%  "

# Row 42
$ws.Range("D42").Value = "'73.19"
$ws.Range("E42").Value = "  +1.14%  "

# Row 43
$ws.Range("D43").Value = "'0.8877"
$ws.Range("E43").Value = "  -1.86%  "

# Row 44
$ws.Range("D44").Value = "'106.63"
$ws.Range("E44").Value = "  +4.85%  "

# Row 45
$ws.Range("D45").Value = "'0.9997"
$ws.Range("E45").Value = "  -0.04%  "

# Row 46
$ws.Range("D46").Value = "2.025.73"
$ws.Range("E46").Value = "  +0.81%  "

# Row 47
$ws.Range("D47").Value = "'1.814"
$ws.Range("E47").Value = "  +2.88%  "

# Row 48
$ws.Range("D48").Value = "'0.5201"
$ws.Range("E48").Value = "  -0.58%  "

# Row 50
$ws.Range("D50").Value = "'9.378"
$ws.Range("E50").Value = "  +0.83%  "

# Row 51
$ws.Range("D51").Value = "'0.4304"
$ws.Range("E51").Value = "  +1.18%  "

Write-Output "Crypto price/volume update applied"
